$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the rich-text "duty cycle" warning in A3 with plain simplified text
$ws.Range("A3").Value = "Time and Date Testing Was Complete (Testing can be completed any time/date):"

# Row 3 no longer needs the taller height that accommodated the old multi-line warning
$ws.Range("A3").EntireRow.RowHeight = 15.75

# Move the active selection to B11
$ws.Range("B11").Select()
